$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.972.77"
$ws.Range("E2").Value = "  +1.20%  "

$ws.Range("D3").Value = "1.894.45"
$ws.Range("E3").Value = "  +0.84%  "

$ws.Range("D4").Value = "1.015"
$ws.Range("E4").Value = "  +1.32%  "

$ws.Range("D5").Value = "335.73"
$ws.Range("E5").Value = "  +1.45%  "

$ws.Range("E6").Value = "  +1.32%  "

$ws.Range("D7").Value = "0.4698"
$ws.Range("E7").Value = "  -0.26%  "

$ws.Range("D8").Value = "0.3934"
$ws.Range("E8").Value = "  -0.79%  "

$ws.Range("E9").Value = "  -0.39%  "

$ws.Range("D10").Value = "0.08066"
$ws.Range("E10").Value = "  +0.11%  "

$ws.Range("D11").Value = "1.023"
$ws.Range("E11").Value = "  -0.22%  "

$ws.Range("D12").Value = "21.86"
$ws.Range("E12").Value = "  +0.37%  "

$ws.Range("D13").Value = "1.899.93"
$ws.Range("E13").Value = "  +0.94%  "

$ws.Range("D14").Value = "5.975"
$ws.Range("E14").Value = "  +0.41%  "

$ws.Range("D15").Value = "7.141"
$ws.Range("E15").Value = "  -0.59%  "

$ws.Range("D16").Value = "1.019"
$ws.Range("E16").Value = "  +1.47%  "

$ws.Range("D17").Value = "0.06802"
$ws.Range("E17").Value = "  +3.38%  "

$ws.Range("D18").Value = "0.00001052"
$ws.Range("E18").Value = "  +1.28%  "

$ws.Range("D19").Value = "87.43"
$ws.Range("E19").Value = "  +0.68%  "

$ws.Range("D20").Value = "17.23"
$ws.Range("E20").Value = "  -0.14%  "

$ws.Range("D21").Value = "1.014"
$ws.Range("E21").Value = "  +1.20%  "

$ws.Range("D22").Value = "28.025.62"
$ws.Range("E22").Value = "  +1.31%  "

$ws.Range("D23").Value = "5.522"
$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("D24").Value = "11.02"
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("D25").Value = "2.348"
$ws.Range("E25").Value = "  +1.68%  "

$ws.Range("D26").Value = "2.113.70"
$ws.Range("E26").Value = "  +0.43%  "

$ws.Range("D27").Value = "159.79"
$ws.Range("E27").Value = "  +3.38%  "

$ws.Range("D28").Value = "20.08"
$ws.Range("E28").Value = "  -0.82%  "

$ws.Range("D29").Value = "2.092"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("D30").Value = "5.483"
$ws.Range("E30").Value = "  -1.86%  "

$ws.Range("D31").Value = "122.02"
$ws.Range("E31").Value = "  -0.30%  "

$ws.Range("D32").Value = "0.9745"
$ws.Range("E32").Value = "  +2.11%  "

$ws.Range("D33").Value = "0.09508"
$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("D34").Value = "3.642"
$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("D35").Value = "1.408"
$ws.Range("E35").Value = "  -4.24%  "

$ws.Range("D36").Value = "5.371"
$ws.Range("E36").Value = "  +1.36%  "

$ws.Range("D37").Value = "0.06151"
$ws.Range("E37").Value = "  +0.71%  "

$ws.Range("D38").Value = "0.02263"
$ws.Range("E38").Value = "  +0.46%  "

$ws.Range("D39").Value = "1.217"
$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("E40").Value = "  -1.12%  "

$ws.Range("D41").Value = "0.6010"
$ws.Range("E41").Value = "  +0.29%  "

$ws.Range("D42").Value = "0.1893"
$ws.Range("E42").Value = "  -0.13%  "

$ws.Range("D43").Value = "10.31"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("E44").Value = "  +1.42%  "

$ws.Range("D45").Value = "0.5715"
$ws.Range("E45").Value = "  +0.31%  "

$ws.Range("D46").Value = "12.21"
$ws.Range("E46").Value = "  +0.52%  "

$ws.Range("B47").Value = "PaxosStandard"
$ws.Range("C47").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D47").Value = "1.133"
$ws.Range("E47").Value = "  +12.70%  "

$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "3.406"
$ws.Range("E48").Value = "  -0.20%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.942"
$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.06932"
$ws.Range("E50").Value = "  +1.59%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "114.04"
$ws.Range("E51").Value = "  +3.68%  "
